# Orden de produccion - agrega el detalle de dietas del Desayuno
# para la fecha 2025-05-15 (serial 45792), manteniendo la fila 2
# existente (antes era la unica fila de datos, con fecha 45782) y
# sumando las filas 3 a 23 con el resto de las dietas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fecha = 45792

$data = @(
    @("Astringente",               "Desayuno", 3),
    @("Blanda",                    "Desayuno", 19),
    @("Coronaria",                 "Desayuno", 15),
    @("Hepatica",                  "Desayuno", 1),
    @("Hipercalorica",             "Desayuno", 1),
    @("Hiperproteica",             "Desayuno", 2),
    @("Hipo Grasa",                "Desayuno", 6),
    @("Hipoglucida",               "Desayuno", 12),
    @("Hiposodica",                "Desayuno", 28),
    @("Liquida Clara",             "Desayuno", 4),
    @("Liquida Total",             "Desayuno", 7),
    @("Liquida Total 140 Cc",      "Desayuno", 1),
    @("Liquida Total Miel 140 Cc", "Desayuno", 2),
    @("Liquida Total Nectar",      "Desayuno", 9),
    @("Liquida Total Nectar 140 Cc","Desayuno", 2),
    @("Liquida total Miel",        "Desayuno", 3),
    @("Normal",                    "Desayuno", 57),
    @("Renal Dialisis",            "Desayuno", 4),
    @("Renal PRE Dialisis",        "Desayuno", 5),
    @("Semiblanda",                "Desayuno", 28),
    @("Semiblanda Pequena",        "Desayuno", 4),
    @("Todo Pure",                 "Desayuno", 2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]

    $ws.Cells.Item($row, 1).Value = $fecha
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 2).Value = $item[0]
    $ws.Cells.Item($row, 3).Value = $item[1]
    $ws.Cells.Item($row, 4).Value = $item[2]
}
